$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append: row 81 (matches format/layout of preceding data rows)

# Copy formatting (styles/number formats) from the last existing data row (80)
# into the new row so A gets the bold/border/centered style (s=1) and E gets
# the date-time number format (s=2), matching every other data row.
$ws.Range("A80:V80").Copy()
$ws.Range("A81:V81").PasteSpecial(-4122)

# Populate the new row's values
$ws.Range("A81").Value = 80
$ws.Range("B81").Value = "croatia"
$ws.Range("C81").Value = "hnl"
$ws.Range("D81").Value = "2023-2024"
$ws.Range("E81").Value = 45262.72916666666
$ws.Range("F81").Value = "Hajduk Split"
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = "Gorica"
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1.33
$ws.Range("K81").Value = "26/11/2023 17:12"
$ws.Range("L81").Value = 1.4
$ws.Range("M81").Value = "02/12/2023 17:28"
$ws.Range("N81").Value = 5.01
$ws.Range("O81").Value = "26/11/2023 17:12"
$ws.Range("P81").Value = 4.34
$ws.Range("Q81").Value = "02/12/2023 17:28"
$ws.Range("R81").Value = 9.12
$ws.Range("S81").Value = "26/11/2023 17:12"
$ws.Range("T81").Value = 9.53
$ws.Range("U81").Value = "02/12/2023 17:28"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/croatia/hnl/hajduk-split-hnk-gorica/0hppQlIp/"
